$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '51.986.04'
$cell.Style = 'Normal'

$cell = $ws.Range('E2')
$cell.NumberFormat = '@'
$cell.Value = '  +0.30%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '3.003.54'
$cell.Style = 'Normal'

$cell = $ws.Range('E3')
$cell.NumberFormat = '@'
$cell.Value = '  +2.54%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E4')
$cell.NumberFormat = '@'
$cell.Value = '  -0.01%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '355.38'
$cell.Style = 'Normal'

$cell = $ws.Range('E5')
$cell.NumberFormat = '@'
$cell.Value = '  -0.09%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '106.69'
$cell.Style = 'Normal'

$cell = $ws.Range('E6')
$cell.NumberFormat = '@'
$cell.Value = '  -3.84%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.555'
$cell.Style = 'Normal'

$cell = $ws.Range('E7')
$cell.NumberFormat = '@'
$cell.Value = '  -2.17%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E8')
$cell.NumberFormat = '@'
$cell.Value = '  +0.19%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.607'
$cell.Style = 'Normal'

$cell = $ws.Range('E9')
$cell.NumberFormat = '@'
$cell.Value = '  -3.42%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '37.97'
$cell.Style = 'Normal'

$cell = $ws.Range('E10')
$cell.NumberFormat = '@'
$cell.Value = '  -3.59%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E11')
$cell.NumberFormat = '@'
$cell.Value = '  +2.74%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '0.0854'
$cell.Style = 'Normal'

$cell = $ws.Range('E12')
$cell.NumberFormat = '@'
$cell.Value = '  -3.19%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E13')
$cell.NumberFormat = '@'
$cell.Value = '  -3.66%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '3.476.61'
$cell.Style = 'Normal'

$cell = $ws.Range('E14')
$cell.NumberFormat = '@'
$cell.Value = '  +2.57%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '7.58'
$cell.Style = 'Normal'

$cell = $ws.Range('E15')
$cell.NumberFormat = '@'
$cell.Value = '  -4.28%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '2.974.98'
$cell.Style = 'Normal'

$cell = $ws.Range('E16')
$cell.NumberFormat = '@'
$cell.Value = '  +1.48%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E17')
$cell.NumberFormat = '@'
$cell.Value = '  +1.55%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '51.987.99'
$cell.Style = 'Normal'

$cell = $ws.Range('E18')
$cell.NumberFormat = '@'
$cell.Value = '  +0.25%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E19')
$cell.NumberFormat = '@'
$cell.Value = '  +2.15%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '7.44'
$cell.Style = 'Normal'

$cell = $ws.Range('E20')
$cell.NumberFormat = '@'
$cell.Value = '  -1.72%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '13.48'
$cell.Style = 'Normal'

$cell = $ws.Range('E21')
$cell.NumberFormat = '@'
$cell.Value = '  -4.07%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E22')
$cell.NumberFormat = '@'
$cell.Value = '  -1.47%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '69.05'
$cell.Style = 'Normal'

$cell = $ws.Range('E23')
$cell.NumberFormat = '@'
$cell.Value = '  -2.70%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '263.69'
$cell.Style = 'Normal'

$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '2.71'
$cell.Style = 'Normal'

$cell = $ws.Range('E25')
$cell.NumberFormat = '@'
$cell.Value = '  -4.08%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E26')
$cell.NumberFormat = '@'
$cell.Value = '  -2.98%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '26.91'
$cell.Style = 'Normal'

$cell = $ws.Range('E27')
$cell.NumberFormat = '@'
$cell.Value = '  -1.08%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E28')
$cell.NumberFormat = '@'
$cell.Value = '  -0.03%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E29')
$cell.NumberFormat = '@'
$cell.Value = '  -0.75%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '0.109'
$cell.Style = 'Normal'

$cell = $ws.Range('E30')
$cell.NumberFormat = '@'
$cell.Value = '  +2.62%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E31')
$cell.NumberFormat = '@'
$cell.Value = '  +4.45%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '10.14'
$cell.Style = 'Normal'

$cell = $ws.Range('E32')
$cell.NumberFormat = '@'
$cell.Value = '  -4.36%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E33')
$cell.NumberFormat = '@'
$cell.Value = '  -8.14%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E34')
$cell.NumberFormat = '@'
$cell.Value = '  +13.19%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '51.26'
$cell.Style = 'Normal'

$cell = $ws.Range('E35')
$cell.NumberFormat = '@'
$cell.Value = '  -1.82%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.0432'
$cell.Style = 'Normal'

$cell = $ws.Range('E36')
$cell.NumberFormat = '@'
$cell.Value = '  -2.64%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E37')
$cell.NumberFormat = '@'
$cell.Value = '  +0.00%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E38')
$cell.NumberFormat = '@'
$cell.Value = '  +0.88%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '2.84'
$cell.Style = 'Normal'

$cell = $ws.Range('E39')
$cell.NumberFormat = '@'
$cell.Value = '  +2.91%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E40')
$cell.NumberFormat = '@'
$cell.Value = '  -3.77%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '17.45'
$cell.Style = 'Normal'

$cell = $ws.Range('E41')
$cell.NumberFormat = '@'
$cell.Value = '  -5.99%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E42')
$cell.NumberFormat = '@'
$cell.Value = '  -2.98%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '23.33'
$cell.Style = 'Normal'

$cell = $ws.Range('E43')
$cell.NumberFormat = '@'
$cell.Value = '  +0.47%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '123.98'
$cell.Style = 'Normal'

$cell = $ws.Range('E44')
$cell.NumberFormat = '@'
$cell.Value = '  +3.83%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E45')
$cell.NumberFormat = '@'
$cell.Value = '  +0.24%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '2.123.24'
$cell.Style = 'Normal'

$cell = $ws.Range('E46')
$cell.NumberFormat = '@'
$cell.Value = '  -0.78%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E47')
$cell.NumberFormat = '@'
$cell.Value = '  -4.60%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '0.243'
$cell.Style = 'Normal'

$cell = $ws.Range('E49')
$cell.NumberFormat = '@'
$cell.Value = '  -2.70%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '0.0331'
$cell.Style = 'Normal'

$cell = $ws.Range('E50')
$cell.NumberFormat = '@'
$cell.Value = '  -1.18%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '0.897'
$cell.Style = 'Normal'

$cell = $ws.Range('E51')
$cell.NumberFormat = '@'
$cell.Value = '  -1.36%  '
$cell.Style = 'Normal'

